$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest scrape.
# A handful of D-column values are plain decimals (e.g. "1.00", "0.999") that
# Excel would otherwise auto-convert to numbers on assignment. Flip those
# specific cells to Text just long enough to assign the literal, then restore
# General formatting so the cell style matches the rest of the (unstyled) sheet.

$ws.Range("D2").Value = "63.906.97"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "2.631.75"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.78"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.57"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.107"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.67"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.50"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "3.098.09"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "63.723.45"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000147"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "2.583.15"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.26"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "343.82"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.94"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.60"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +10.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.65"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.90"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.39"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "500.09"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +8.76%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.99"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +3.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.76"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +9.57%  "
$ws.Range("D33").Value = "0.0₃0814"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "174.96"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.399"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.05"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.75"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "165.66"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.35"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.76"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.69"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +5.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.630"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0545"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.64"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.73"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  -0.38%  "
